# Added fromMeeting and fromOffer text checkers
# Update the sample/example row with new reference values.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("A2").Value = "ME-562"
$ws.Range("B2").Value = "OF-567"
